# Apply "Actualización automática 2025-05-30 16:20:08" changes to the
# "VENTAS POR GRUPO" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

# --- Column L (PORCELANATO) and other per-row updates ---

# Row 2: BARRAGAN PUENTE NATALY CAROLINA
$ws.Range("L2").Value = -49.25

# Row 3: BECERRA FARIAS ROSA DAYANA
$ws.Range("L3").Value = 1451.52

# Row 5: CAIZA COLLAGUAZO ROCIO PILAR
$ws.Range("L5").Value = 226.8

# Row 6: CERAMICAS AL COSTO S.A.S.
$ws.Range("E6").Value = 208.35
$ws.Range("F6").Value = 166.32
$ws.Range("L6").Value = 2492.93

# Row 8: CONZA VEGA FRANCO BLADYMIR
$ws.Range("L8").Value = 2785.1

# Row 10: F.V - AREA ANDINA S.A.
$ws.Range("L10").Value = 156.67

# Row 14: TAMAYO VILLACIS EDWIN XAVIER
$ws.Range("L14").Value = 30.61

# Row 15: TOSCANO RAMIREZ MONICA CECILIA
$ws.Range("L15").Value = 4277.03

# Row 18: ZAMBRANO ANGELA MARIA
$ws.Range("D18").Value = 10277.38
$ws.Range("L18").Value = 13819.55

# --- Row 19 summary counts ---
$ws.Range("E19").Value = "1 de 17"
$ws.Range("F19").Value = "1 de 17"
$ws.Range("L19").Value = "8 de 17"
